# Update "PatternRecognition-Schedule" worksheet:
#  - mark Progress (column F) as 90% complete for sessions 1-14 (rows 3-16)
#  - add a new session 15 (row 17): Progress 60%, Subject "PCA",
#    Problem "didn't get clear on it"
#  - move the on-screen selection to H17 (the newly-filled cell)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Progress column (F) for rows 3 through 16 -> 90%
foreach ($r in 3..16) {
    $ws.Cells.Item($r, 6).Value = 0.9
}

# Row 17: new "PCA" session
$ws.Range("F17").Value = 0.6
$ws.Range("G17").Value = "PCA"
$ws.Range("H17").Value = "didn't get clear on it"

# Reflect the edit location in the saved view
$ws.Range("H17").Select()
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 5
    $win.ScrollColumn = 1
} catch {
    # ScrollRow/ScrollColumn not critical to the data edit; ignore if unsupported
}
